# MSME Country Indicators - Kuwait: "Data" -> "Summary" sheet update.
# Renames the sheet and re-lays-out the MSME participation block:
#   - adds a new "Source Type" sub-heading (bold + underlined) at row 9
#   - moves the Micro / SMEs / MSMEs column headers from row 5 to row 11
#   - moves the "Enterprises (% of total)" label from row 6 to row 12,
#     and adds the 90 data point next to it in D12
#   - adds a new italic source citation at row 13

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet used to be called "Data"; it is now the "Summary" sheet.
$ws.Name = "Summary"

# The old header row (5) and label row (6) are vacated - their content
# is re-created further down the sheet (rows 11/12) to make room for the
# new "Source Type" sub-heading above them.
$ws.Range("B5:D5").Clear()
$ws.Range("A6").Clear()

# New bold + underlined sub-heading.
$ws.Range("A9").Value = "Source Type: SME Associations"
$ws.Range("A9").Font.Bold = $true
$ws.Range("A9").Font.Underline = $true

# Column headers (bold "title" style), re-inserted three rows lower.
$ws.Range("B11").Value = "Micro"
$ws.Range("B11").Font.Bold = $true
$ws.Range("C11").Value = "SMEs"
$ws.Range("C11").Font.Bold = $true
$ws.Range("D11").Value = "MSMEs"
$ws.Range("D11").Font.Bold = $true

# Data label (bold) plus its value, stored as text.
$ws.Range("A12").Value = "Enterprises (% of total)"
$ws.Range("A12").Font.Bold = $true
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "90"

# New italic source citation.
$ws.Range("A13").Value = "Source: WAMDA, 2007"
$ws.Range("A13").Font.Italic = $true
